$wb = $excel.ActiveWorkbook

function Add-TradeRow {
    $ws = $global:targetSheet
    $row = 78

    $ws.Cells.Item($row, 1).Value = 77

    # Date column - force text storage so "2026-02-17" is not coerced into a date serial value,
    # then clear the formatting we applied so no extra style is retained.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.ClearFormats()

    $ws.Cells.Item($row, 3).Value = "12:57:47"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.51

    # Exit Price - trade is still OPEN, so this column is blank.
    $ws.Cells.Item($row, 7).Font.Italic = $false

    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100.1315249294667
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"

    # Exit Reason - trade is still OPEN, so this column is blank.
    $ws.Cells.Item($row, 16).Font.Italic = $false

    $ws.Cells.Item($row, 17).Value = 0
}

$global:targetSheet = $wb.Worksheets.Item("All Trades")
Add-TradeRow

$global:targetSheet = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow
